$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44701, 200, 'Región del Maule'),
    @(44726, 150, 'Región del Maule'),
    @(44354, 150, 'Región del Maule'),
    @(44372, 150, 'Región Metropolitana'),
    @(44720, 150, 'Región Metropolitana'),
    @(44355, 150, 'Región Metropolitana'),
    @(44371, 150, 'Región Metropolitana'),
    @(44719, 150, 'Región Metropolitana'),
    @(44342, 150, 'Región del Maule'),
    @(44348, 150, 'Región del Maule'),
    @(44707, 150, 'Región Metropolitana'),
    @(44386, 200, 'Región Metropolitana'),
    @(44376, 150, 'Región Metropolitana'),
    @(44340, 150, 'Región del Maule'),
    @(44715, 150, 'Región Metropolitana'),
    @(44362, 100, 'Región Metropolitana'),
    @(44690, 500, 'Región del Maule'),
    @(44725, 150, 'Región del Maule'),
    @(44711, 150, 'Región Metropolitana'),
    @(44706, 200, 'Región Metropolitana'),
    @(44358, 150, 'Región Metropolitana'),
    @(44364, 100, 'Región Metropolitana')
)

$prices = @(
    @(7000, 7000, 7000),
    @(8000, 8000, 8000),
    @(7000, 7000, 7000),
    @(7000, 7000, 7000),
    @(9000, 9000, 9000),
    @(7000, 7000, 7000),
    @(6500, 6500, 6500),
    @(9000, 9000, 9000),
    @(7000, 7000, 7000),
    @(7000, 7000, 7000),
    @(9000, 9000, 9000),
    @(6500, 6500, 6500),
    @(6500, 6500, 6500),
    @(7000, 7000, 7000),
    @(9000, 9000, 9000),
    @(6500, 6500, 6500),
    @(7000, 7000, 7000),
    @(8000, 8000, 8000),
    @(8500, 8500, 8500),
    @(9000, 9000, 9000),
    @(7000, 7000, 7000),
    @(7000, 7000, 7000)
)

$pkg = @(194, 222, 194, 194, 250, 194, 181, 250, 194, 194, 250, 181, 181, 194, 250, 181, 194, 222, 236, 250, 194, 194)

for ($i = 0; $i -lt 22; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 4).Value = $data[$i][0]
    $ws.Cells.Item($r, 10).Value = $data[$i][1]
    $ws.Cells.Item($r, 11).Value = $prices[$i][0]
    $ws.Cells.Item($r, 12).Value = $prices[$i][1]
    $ws.Cells.Item($r, 13).Value = $prices[$i][2]
    $ws.Cells.Item($r, 15).Value = $data[$i][2]
    $ws.Cells.Item($r, 16).Value = $pkg[$i]
}
